# Script 1 - atualização automática de dados (Execução: 22)
# Refreshes the "g1.1" sheet data: UF rankings move from the 2022 vintage
# to the 2023 vintage (values, associated ranks, UF order and the
# "Variação (%) 20XX" / "Variação (%) 20XX/2010" category labels).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: rows 2-10 -> "Variação (%) 2023" -------------------------
$block1 = @(
    @{ Row = 2;  UF = "AC"; Valor = 14.73257689442189; Rank = 1 },
    @{ Row = 3;  UF = "MS"; Valor = 13.44269577606423; Rank = 2 },
    @{ Row = 4;  UF = "MT"; Valor = 12.88001598426398; Rank = 3 },
    @{ Row = 5;  UF = "TO"; Valor = 7.890383025089162; Rank = 4 },
    @{ Row = 6;  UF = "RJ"; Valor = 5.652659822157795; Rank = 5 },
    @{ Row = 7;  UF = "GO"; Valor = 4.816953216278661; Rank = 6 },
    @{ Row = 8;  UF = "SE"; Valor = 3.118144130554446; Rank = 15 },
    @{ Row = 9;  UF = "BR"; Valor = 3.241657824791806; Rank = $null },
    @{ Row = 10; UF = "NE"; Valor = 2.867008788862638; Rank = $null }
)

foreach ($item in $block1) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.UF
    $ws.Cells.Item($r, 2).Value = $item.Valor
    if ($item.Rank -ne $null) {
        $ws.Cells.Item($r, 3).Value = $item.Rank
    }
    $ws.Cells.Item($r, 4).Value = "Variação (%) 2023"
}

# --- Block 2: rows 11-19 -> "Variação (%) 2023/2010" --------------------
$block2 = @(
    @{ Row = 11; UF = "MT"; Valor = 77.14346626765018; Rank = 1 },
    @{ Row = 12; UF = "RR"; Valor = 63.84407511155798; Rank = 2 },
    @{ Row = 13; UF = "TO"; Valor = 59.35173933449352; Rank = 3 },
    @{ Row = 14; UF = "MS"; Valor = 49.57179111911111; Rank = 4 },
    @{ Row = 15; UF = "AC"; Valor = 42.69524774665621; Rank = 5 },
    @{ Row = 16; UF = "PI"; Valor = 36.70681689547283; Rank = 6 },
    @{ Row = 17; UF = "SE"; Valor = 7.324239245718005; Rank = 25 },
    @{ Row = 18; UF = "BR"; Valor = 15.14430956101356; Rank = $null },
    @{ Row = 19; UF = "NE"; Valor = 16.86384673819174; Rank = $null }
)

foreach ($item in $block2) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.UF
    $ws.Cells.Item($r, 2).Value = $item.Valor
    if ($item.Rank -ne $null) {
        $ws.Cells.Item($r, 3).Value = $item.Rank
    }
    $ws.Cells.Item($r, 4).Value = "Variação (%) 2023/2010"
}

Write-Output "g1.1 data refreshed to 2023 vintage"
